$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 371 (shifts existing rows 371:407 down to 372:408,
# copying formatting from the row above - matches the date style on column D).
$ws.Rows("371:371").Insert()

# Populate the newly inserted row 371 with the new data record.
$ws.Range("A371").Value = 10
$ws.Range("B371").Value = "Vega Modelo de Temuco"
$ws.Range("C371").Value = "La Araucanía"
$ws.Range("D371").Value = 44578
$ws.Range("E371").Value = 9
$ws.Range("F371").Value = 100112043
$ws.Range("G371").Value = "Pepino ensalada"
$ws.Range("H371").Value = "Sin especificar"
$ws.Range("I371").Value = "Primera"
$ws.Range("J371").Value = 600
$ws.Range("K371").Value = 12000
$ws.Range("L371").Value = 12000
$ws.Range("M371").Value = 12000
$ws.Range("N371").Value = "$/caja 60 unidades"
$ws.Range("O371").Value = "Región del Maule"
$ws.Range("P371").Value = 200
$ws.Range("Q371").Value = 60
$ws.Range("R371").Value = "Hortaliza"
